# Laravel Form & CRUD - insert new "Create Function" slide
# right after "Index Function" (slide 13) and before "Store Function"
# (which was slide 14 and becomes slide 15).

$p = $ppt.ActivePresentation

$nbsp = [char]0x00A0
$lsq  = [char]0x2018   # left single quotation mark
$rsq  = [char]0x2019   # right single quotation mark
$tab  = [char]0x0009

# Insert a new slide at position 14 using the same "Title and Content"
# layout (ppLayoutText = 16) the surrounding slides use.
$null = $p.Slides.Add(14, 16)
$s = $p.Slides.Item(14)

# ---- Title -------------------------------------------------------------
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Create Function"

# ---- Body ----------------------------------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange

$line1 = "Edit the create function in StudentController.php" + $nbsp + "file under" + $nbsp + "app/Http/Controllers" + $nbsp + "directory"
$line2 = ""
$line3 = "public function create()"
$line4 = "{"
$line5 = $tab + " return view(" + $lsq + "create" + $rsq + ");"
$line6 = "}"

$body.Text = $line1 + "`r" + $line2 + "`r" + $line3 + "`r" + $line4 + "`r" + $line5 + "`r" + $line6

# bodyPr autofit: full size (no fontScale/lnSpcReduction)
$body.Parent.WordWrap = $true

# ---- Paragraph 1 runs: bold "StudentController.php" + trailing nbsp ------
$para1 = $body.Paragraphs(1)
$run1bStart = "Edit the create function in ".Length + 1
$run1bLen   = "StudentController.php".Length
$para1.Characters($run1bStart, $run1bLen).Font.Bold = $true
$para1.Characters($run1bStart + $run1bLen, 1).Font.Bold = $true

# ---- Paragraphs 3-6: Courier New "code" lines, no bullet -----------------
for ($i = 3; $i -le 6; $i++) {
    $para = $body.Paragraphs($i)
    $para.ParagraphFormat.Bullet.Visible = 0
    $allRun = $para.Characters(1, $para.Text.Length)
    $allRun.Font.Bold = $true
    $allRun.Font.Name = "Courier New"
    $allRun.Font.Size = 22
}

# Paragraph 5 ("\t return view('create');") splits into two runs:
# the leading tab stays at size 22, the rest is size 24.
$para5 = $body.Paragraphs(5)
$para5.Characters(2, $para5.Text.Length - 1).Font.Size = 24

Write-Output "Inserted Create Function slide at index 14; total slides: $($p.Slides.Count)"
Write-Output ("Slide order: " + ((1..$p.Slides.Count | ForEach-Object { $p.Slides.Item($_).Shapes.Item(1).TextFrame.TextRange.Text }) -join " | "))
